# "Actualizar" automatic update: the "Fecha" (Date) column (D) holds a
# rolling log of update timestamps. Each run of the updater:
#   - stamps the most-recent block of rows (2-15) with the new "now" time,
#   - and shifts every older block down one slot (16-29 <- old 2-15,
#     30-43 <- old 16-29), dropping the oldest block off the bottom.
#
# Concretely (Excel serial date-times):
#   rows 30-43 <- 44302.46751424768   (was the old rows 16-29 value)
#   rows 16-29 <- 44302.48887386574   (was ~ the old rows 2-15 value)
#   rows 2-15  <- 44302.51026142319   (brand new update timestamp)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldestBlock = 44302.46751424768
$middleBlock = 44302.48887386574
$newestBlock = 44302.51026142319

# Shift the oldest surviving block down into rows 30-43.
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldestBlock
}

# Shift the previous newest block down into rows 16-29.
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $middleBlock
}

# Stamp rows 2-15 with the freshly captured update timestamp.
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newestBlock
}
